$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:58:43 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:59:45 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 21:00:52 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 21:01:58 EST 2024"
$ws.Cells.Item(6, 2).Value = "Tue Dec 03 21:03:03 EST 2024"
$ws.Cells.Item(7, 2).Value = "Tue Dec 03 21:04:04 EST 2024"
$ws.Cells.Item(8, 2).Value = "Tue Dec 03 21:05:11 EST 2024"
$ws.Cells.Item(9, 2).Value = "Tue Dec 03 21:06:17 EST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:53:31 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:54:55 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:56:14 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:58:40 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowCreditSCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:06:29 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:07:48 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:09:12 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:10:32 EST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckSCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:27:42 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:29:08 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:30:30 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:32:57 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowPersonalSavingsSCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:42:58 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:44:24 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:45:45 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:48:13 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowCreditDCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:01:02 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:02:21 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:03:45 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:05:04 EST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:37:55 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:39:20 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:40:40 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:43:07 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowPC")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:11:57 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:13:05 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:14:07 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:15:14 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"
$ws.Cells.Item(6, 2).Value = "Tue Dec 03 20:16:53 EST 2024"
$ws.Cells.Item(6, 1).Value = "Fail"
$ws.Cells.Item(7, 2).Value = "Tue Dec 03 20:18:27 EST 2024"
$ws.Cells.Item(7, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckDCF")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:35:20 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:36:47 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:38:08 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:40:34 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowPS")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 20:50:36 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 20:51:44 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 20:52:46 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 20:53:52 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"
$ws.Cells.Item(6, 2).Value = "Tue Dec 03 20:55:31 EST 2024"
$ws.Cells.Item(6, 1).Value = "Fail"
$ws.Cells.Item(7, 2).Value = "Tue Dec 03 20:57:04 EST 2024"
$ws.Cells.Item(7, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("OverAndUnderPayCredit")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:22:42 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:23:47 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:24:55 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:26:02 EST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayPC")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:27:07 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:28:12 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:29:22 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:30:58 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("OverAndUnderPayPS")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:32:34 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:33:39 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:34:44 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:36:19 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("OverAndUnderPayCorp")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:17:24 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:18:27 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:19:30 EST 2024"
$ws.Cells.Item(4, 1).Value = "Fail"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:21:05 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("PayNowCorp")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 19:45:29 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 19:46:36 EST 2024"
$ws.Cells.Item(4, 2).Value = "Tue Dec 03 19:47:38 EST 2024"
$ws.Cells.Item(5, 2).Value = "Tue Dec 03 19:48:42 EST 2024"
$ws.Cells.Item(5, 1).Value = "Fail"
$ws.Cells.Item(6, 2).Value = "Tue Dec 03 19:50:20 EST 2024"
$ws.Cells.Item(6, 1).Value = "Fail"
$ws.Cells.Item(7, 2).Value = "Tue Dec 03 19:51:54 EST 2024"
$ws.Cells.Item(7, 1).Value = "Fail"

$ws = $wb.Worksheets.Item("CardNotAcceptedErrorCC")
$ws.Cells.Item(2, 2).Value = "Tue Dec 03 21:07:23 EST 2024"
$ws.Cells.Item(3, 2).Value = "Tue Dec 03 21:08:29 EST 2024"
